$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'99.134.07"
$ws.Range("E2").Value = "'  +1.57%  "

$ws.Range("D3").Value = "'3.309.28"
$ws.Range("E3").Value = "'  -0.90%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'255.71"
$ws.Range("E5").Value = "'  -0.73%  "

$ws.Range("D6").Value = "'623.55"
$ws.Range("E6").Value = "'  +0.58%  "

$ws.Range("D7").Value = "'1.47"
$ws.Range("E7").Value = "'  +30.93%  "

$ws.Range("D8").Value = "'0.409"
$ws.Range("E8").Value = "'  +5.98%  "

$ws.Range("E9").Value = "'  +0.00%  "

$ws.Range("D10").Value = "'0.976"
$ws.Range("E10").Value = "'  +22.79%  "

$ws.Range("D11").Value = "'3.306.23"
$ws.Range("E11").Value = "'  -0.96%  "

$ws.Range("E12").Value = "'  +0.43%  "

$ws.Range("D13").Value = "'39.62"
$ws.Range("E13").Value = "'  +11.46%  "

$ws.Range("D14").Value = "'98.936.54"
$ws.Range("E14").Value = "'  +1.75%  "

$ws.Range("E15").Value = "'  +1.70%  "

$ws.Range("D16").Value = "'3.931.88"
$ws.Range("E16").Value = "'  -0.63%  "

$ws.Range("E17").Value = "'  -0.65%  "

$ws.Range("D18").Value = "'3.309.40"
$ws.Range("E18").Value = "'  -0.71%  "

$ws.Range("D19").Value = "'3.47"
$ws.Range("E19").Value = "'  -3.11%  "

$ws.Range("D20").Value = "'15.59"
$ws.Range("E20").Value = "'  +4.09%  "

$ws.Range("D21").Value = "'6.30"
$ws.Range("E21").Value = "'  +8.26%  "

$ws.Range("D22").Value = "'486.89"
$ws.Range("E22").Value = "'  +0.86%  "

$ws.Range("D23").Value = "'9.46"
$ws.Range("E23").Value = "'  +2.44%  "

$ws.Range("E24").Value = "'  -2.19%  "

$ws.Range("D25").Value = "'5.64"
$ws.Range("E25").Value = "'  +0.07%  "

$ws.Range("D26").Value = "'89.06"
$ws.Range("E26").Value = "'  +1.18%  "

$ws.Range("D27").Value = "'12.00"
$ws.Range("E27").Value = "'  -0.93%  "

$ws.Range("D28").Value = "'0.306"
$ws.Range("E28").Value = "'  +28.20%  "

$ws.Range("D29").Value = "'3.490.33"
$ws.Range("E29").Value = "'  -0.90%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  -0.12%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.135"
$ws.Range("E31").Value = "'  +11.12%  "

$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").Value = "'0.188"
$ws.Range("E32").Value = "'  +2.47%  "

$ws.Range("D33").Value = "'10.29"
$ws.Range("E33").Value = "'  +11.28%  "

$ws.Range("E34").Value = "'  +0.04%  "

$ws.Range("D35").Value = "'27.83"
$ws.Range("E35").Value = "'  +1.58%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'7.23"
$ws.Range("E36").Value = "'  -2.64%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.471"
$ws.Range("E37").Value = "'  +4.67%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.148"
$ws.Range("E38").Value = "'  -1.88%  "

$ws.Range("E39").Value = "'  +0.35%  "

$ws.Range("D40").Value = "'24.84"
$ws.Range("E40").Value = "'  +0.10%  "

$ws.Range("D41").Value = "'490.77"
$ws.Range("E41").Value = "'  -4.03%  "

$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = "'  +3.34%  "

$ws.Range("D43").Value = "'1.23"
$ws.Range("E43").Value = "'  -3.29%  "

$ws.Range("D44").Value = "'0.787"
$ws.Range("E44").Value = "'  -0.70%  "

$ws.Range("E45").Value = "'  +0.00%  "

$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "'  -5.83%  "

$ws.Range("D47").Value = "'1.96"
$ws.Range("E47").Value = "'  +2.02%  "

$ws.Range("D48").Value = "'158.40"
$ws.Range("E48").Value = "'  -1.59%  "

$ws.Range("D49").Value = "'7.32"
$ws.Range("E49").Value = "'  +15.83%  "

$ws.Range("D50").Value = "'0.848"
$ws.Range("E50").Value = "'  +5.46%  "

$ws.Range("D51").Value = "'4.71"
$ws.Range("E51").Value = "'  +4.50%  "
